$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Give row 52 a thin top + thin bottom border (new border/style added to styles.xml)
$rng52 = $ws.Range("A52:E52")
$rng52.Borders(8).Weight = 2
$rng52.Borders(9).Weight = 2

# 2) Add the two new rows of data (53 and 54).
#    Cell values are set in the same column-major order the original author
#    used (C col, then D col, then E col, then the A54 filename, then the
#    remaining A/B cells) so new shared-string entries land in the same order.
$ws.Range("C53").Value = " It's terrible! I don't do anything\nother than trying to win at Big Treasure!"
$ws.Range("C54").Value = " I can't go on like this!"
$ws.Range("D53").Value = " Какой ужас! Я пытаюсь выиграть\nБольшое Сокровище и больше ничего не\nделаю!"
$ws.Range("D54").Value = " Я так больше не могу!"
$ws.Range("E53").Value = " Ëàëïê ôçàò! Ÿ ðúóàýòû âúéãñàóû\nÁïìûšïå Òïëñïâéþå é áïìûšå îéœåãï îå\näåìàý!"
$ws.Range("E54").Value = " Ÿ óàë áïìûšå îå íïãô!"
$ws.Range("A54").Value = "SCRIPT/P01P04A/us2302.ssb"
$ws.Range("A53").Value = "SCRIPT/P01P04A/us2202.ssb"
$ws.Range("B53").Value = 18
$ws.Range("B54").Value = 21

# 3) Match the row heights used by the rest of the sheet for 3-line wrapped text.
$ws.Rows(53).RowHeight = 43.2
$ws.Rows(54).RowHeight = 43.2

# 4) Update the active selection to the new last cell, like the saved workbook.
$ws.Range("E54").Select()
